# Update "想去人数" (column F) counts that were refreshed by the gh-pages
# data generation run (commit 456a3b4). The same underlying event rows are
# duplicated on the "展览" sheet and the "全部类型" sheet, and F32 happens to
# have diverged slightly between the two copies, so each sheet is updated
# with its own (independent) set of new values.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Row -> new value for the "展览" sheet
$sheet1Updates = @{
    3  = 198
    6  = 183
    7  = 258
    8  = 43
    9  = 252
    10 = 15
    11 = 41
    14 = 1529
    15 = 48
    16 = 489
    17 = 446
    19 = 69
    20 = 36
    22 = 1385
    23 = 3327
    27 = 1078
    28 = 77
    29 = 1670
    32 = 50
    33 = 275
    34 = 398
    36 = 635
    38 = 30
}

# Row -> new value for the "全部类型" sheet
$sheet4Updates = @{
    3  = 198
    6  = 183
    7  = 258
    8  = 43
    9  = 252
    10 = 15
    11 = 41
    14 = 1529
    15 = 48
    16 = 489
    17 = 446
    19 = 69
    20 = 36
    22 = 1385
    23 = 3327
    27 = 1078
    28 = 77
    29 = 1670
    32 = 51
    33 = 275
    34 = 398
    36 = 635
    38 = 30
}

foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
